# "Generate Report for Handback"
#
# The two files that were out for localization have now been handed back.
# This updates the localization-status report to reflect that:
#   - The Overview sheet status for zh-cn / de-de changes from
#     "In Translation" to "Handed back: in sync with en-US".
#   - The zh-cn and de-de detail sheets get their "Latest Target File"
#     (hyperlinked .md), "Latest Handback File" and "Latest Handback
#     DateTime" columns populated for both rows, and their "Status"
#     column updated to match the Overview sheet.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$mdFile1 = "06b24daf-f94f-4b72-acdc-6b17a6ea547b.md"
$mdFile2 = "6ec48260-dc32-4c02-b5bd-dd23585f1ec4.md"

$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f252926ebdb1b6400b8360d54be1a3ec8a92e392/e2e/06b24daf-f94f-4b72-acdc-6b17a6ea547b.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f252926ebdb1b6400b8360d54be1a3ec8a92e392/e2e/6ec48260-dc32-4c02-b5bd-dd23585f1ec4.md"

# visual formatting that matches the workbook's existing "HyperLink" style
# (underline + cornflower blue font, same as column A's hyperlinks)
$hyperlinkUnderline = 2
$hyperlinkColor = 15570276

# ---------------------------------------------------------------------
# Overview sheet: update the per-language status cells
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value2 = $newStatus
$overview.Range("F2").Value2 = $newStatus
$overview.Range("E3").Value2 = $newStatus
$overview.Range("F3").Value2 = $newStatus

$overview.Columns.Item(5).AutoFit() | Out-Null
$overview.Columns.Item(6).AutoFit() | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

# Status column (column C) for both rows
$zhcn.Range("C2").Value2 = $newStatus
$zhcn.Range("C3").Value2 = $newStatus

# Row 2 (06b24daf file): Latest Target File (I), Latest Handback File (J),
# Latest Handback DateTime (K)
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrl1, "", "", $mdFile1) | Out-Null
$zhcn.Range("I2").Font.Underline = $hyperlinkUnderline
$zhcn.Range("I2").Font.Color = $hyperlinkColor
$zhcn.Range("J2").Value2 = "06b24daf-f94f-4b72-acdc-6b17a6ea547b.50ff7f55974e55a0adf9f95695e2147967716aba.zh-cn.xlf"
$zhcn.Range("K2").Value2 = "2016-08-24 10:24:25"

# Row 3 (6ec48260 file): Latest Target File (I), Latest Handback File (J),
# Latest Handback DateTime (K)
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $mdUrl2, "", "", $mdFile2) | Out-Null
$zhcn.Range("I3").Font.Underline = $hyperlinkUnderline
$zhcn.Range("I3").Font.Color = $hyperlinkColor
$zhcn.Range("J3").Value2 = "6ec48260-dc32-4c02-b5bd-dd23585f1ec4.bc04f6ddef7ec1125d8d2cb0133eaf2a0eff8de1.zh-cn.xlf"
$zhcn.Range("K3").Value2 = "2016-08-24 10:24:25"

$zhcn.Columns.Item(3).AutoFit() | Out-Null
$zhcn.Columns.Item(9).ColumnWidth = 40
$zhcn.Columns.Item(10).ColumnWidth = 40

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

# Status column (column C) for both rows
$dede.Range("C2").Value2 = $newStatus
$dede.Range("C3").Value2 = $newStatus

# Row 2 (06b24daf file): Latest Target File (I), Latest Handback File (J),
# Latest Handback DateTime (K)
$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrl1, "", "", $mdFile1) | Out-Null
$dede.Range("I2").Font.Underline = $hyperlinkUnderline
$dede.Range("I2").Font.Color = $hyperlinkColor
$dede.Range("J2").Value2 = "06b24daf-f94f-4b72-acdc-6b17a6ea547b.50ff7f55974e55a0adf9f95695e2147967716aba.de-de.xlf"
$dede.Range("K2").Value2 = "2016-08-24 10:24:31"

# Row 3 (6ec48260 file): Latest Target File (I), Latest Handback File (J),
# Latest Handback DateTime (K)
$dede.Hyperlinks.Add($dede.Range("I3"), $mdUrl2, "", "", $mdFile2) | Out-Null
$dede.Range("I3").Font.Underline = $hyperlinkUnderline
$dede.Range("I3").Font.Color = $hyperlinkColor
$dede.Range("J3").Value2 = "6ec48260-dc32-4c02-b5bd-dd23585f1ec4.bc04f6ddef7ec1125d8d2cb0133eaf2a0eff8de1.de-de.xlf"
$dede.Range("K3").Value2 = "2016-08-24 10:24:31"

$dede.Columns.Item(3).AutoFit() | Out-Null
$dede.Columns.Item(9).ColumnWidth = 40
$dede.Columns.Item(10).ColumnWidth = 40
